$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "25.957.83"
Set-TextCell "E2" "  +0.22%  "
Set-TextCell "D3" "1.644.11"
Set-TextCell "E3" "  +0.30%  "
Set-TextCell "D4" "1.008"
Set-TextCell "E4" "  +0.48%  "
Set-TextCell "D5" "215.54"
Set-TextCell "E5" "  -0.16%  "
Set-TextCell "D6" "0.5087"
Set-TextCell "E6" "  +1.03%  "
Set-TextCell "D7" "1.005"
Set-TextCell "E7" "  +0.11%  "
Set-TextCell "D8" "0.2569"
Set-TextCell "E8" "  +0.10%  "
Set-TextCell "D9" "0.06390"
Set-TextCell "E9" "  -0.04%  "
Set-TextCell "D10" "19.53"
Set-TextCell "E10" "  -0.86%  "
Set-TextCell "D11" "0.07794"
Set-TextCell "E11" "  +0.79%  "
Set-TextCell "D12" "4.310"
Set-TextCell "E12" "  +0.98%  "
Set-TextCell "D13" "1.646.94"
Set-TextCell "E13" "  +0.29%  "
Set-TextCell "D14" "0.5474"
Set-TextCell "E14" "  +0.29%  "
Set-TextCell "D15" "0.0₅7871"
Set-TextCell "E15" "  -0.60%  "
Set-TextCell "D16" "64.54"
Set-TextCell "E16" "  +0.10%  "
Set-TextCell "D17" "26.010.20"
Set-TextCell "E17" "  +0.35%  "
Set-TextCell "D18" "1.006"
Set-TextCell "E18" "  +0.21%  "
Set-TextCell "D19" "198.32"
Set-TextCell "E19" "  -2.34%  "
Set-TextCell "D20" "4.435"
Set-TextCell "E20" "  +1.18%  "
Set-TextCell "D21" "9.969"
Set-TextCell "E21" "  +0.53%  "
Set-TextCell "D22" "6.059"
Set-TextCell "E22" "  +1.08%  "
Set-TextCell "D23" "1.005"
Set-TextCell "E23" "  +0.01%  "
Set-TextCell "D24" "1.884"
Set-TextCell "E24" "  -2.45%  "
Set-TextCell "D25" "141.63"
Set-TextCell "E25" "  +0.06%  "
Set-TextCell "E26" "  +0.58%  "
Set-TextCell "D27" "6.888"
Set-TextCell "E27" "  +2.29%  "
Set-TextCell "D28" "15.76"
Set-TextCell "E28" "  +0.29%  "
Set-TextCell "D29" "0.05049"
Set-TextCell "E29" "  +2.39%  "
Set-TextCell "D30" "1.238"
Set-TextCell "E30" "  -0.65%  "
Set-TextCell "D31" "3.266"
Set-TextCell "E31" "  -0.28%  "
Set-TextCell "D32" "3.197"
Set-TextCell "E32" "  +0.22%  "
Set-TextCell "D33" "1.542"
Set-TextCell "E33" "  -0.08%  "
Set-TextCell "D34" "2.363"
Set-TextCell "E34" "  -0.60%  "
Set-TextCell "D35" "0.8962"
Set-TextCell "E35" "  +0.08%  "
Set-TextCell "D36" "2.603"
Set-TextCell "E36" "  -1.14%  "
Set-TextCell "D37" "1.133.42"
Set-TextCell "E37" "  -2.68%  "
Set-TextCell "D38" "0.5497"
Set-TextCell "E38" "  -2.08%  "
Set-TextCell "B39" "VeChain"
Set-TextCell "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D39" "0.01559"
Set-TextCell "E39" "  -0.45%  "
Set-TextCell "B40" "BabyDogeCoin"
Set-TextCell "C40" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D40" "0.0₈133"
Set-TextCell "E40" "  +14.32%  "
Set-TextCell "E41" "  -0.10%  "
Set-TextCell "D42" "2.542"
Set-TextCell "E42" "  -0.97%  "
Set-TextCell "D43" "5.629"
Set-TextCell "E43" "  -1.75%  "
Set-TextCell "D44" "0.8177"
Set-TextCell "E44" "  +1.09%  "
Set-TextCell "D45" "100.01"
Set-TextCell "D46" "1.778.77"
Set-TextCell "E46" "  +0.18%  "
Set-TextCell "D47" "0.4548"
Set-TextCell "E47" "  +0.52%  "
Set-TextCell "D48" "1.003"
Set-TextCell "E48" "  -0.25%  "
Set-TextCell "D49" "54.95"
Set-TextCell "E49" "  -0.06%  "
Set-TextCell "D50" "0.05082"
Set-TextCell "E50" "  +0.49%  "
Set-TextCell "E51" "  +0.16%  "
